{"js": "// Replace long run-on paragraphs with explicit line-break (<w:br/>) separated\n// segments, matching the target OOXML diff. A Word manual line break is\n// represented in the Office.js text model by the vertical-tab character\n// (\\u000b / Chr(11)); inserting that character causes the host to emit a\n// <w:br/> element and split the surrounding text into separate <w:t> runs.\nconst LB = \"\\u000b\";\n\nasync function insertBreaksAfter(paragraph, anchors, breaksPerAnchor) {\n  // anchors: ordered list of unique text snippets; a line break (or several)\n  // is inserted immediately after each snippet is found.\n  for (const anchor of anchors) {\n    const results = paragraph.search(anchor, { matchCase: true });\n    await context.sync();\n    if (results.items.length !== 1) {\n      throw new Error(\n        \"expected exactly 1 match for \" + JSON.stringify(anchor) +\n        \" but found \" + results.items.length\n      );\n    }\n    const found = results.items[0];\n    found.insertText(LB.repeat(breaksPerAnchor), Word.InsertLocation.after);\n    await context.sync();\n  }\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the affected paragraphs by matching their (still joined) text so the\n// script does not depend on brittle paragraph indices.\nfunction findParagraph(substring) {\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.indexOf(substring) !== -1) {\n      return paragraphs.items[i];\n    }\n  }\n  throw new Error(\"paragraph containing \" + JSON.stringify(substring) + \" not found\");\n}\n\n// 1) \"Gerais - ...\" (Portuguese objectives) \u2014 split into two sentences with a\n//    blank line (two <w:br/>) between them.\nconst pGeraisPt = findParagraph(\"Gerais - Apresentar e Ensinar conceitos\");\nawait insertBreaksAfter(\n  pGeraisPt,\n  [\"realizarem as mudan\u00e7as que se fa\u00e7am necess\u00e1rias.\"],\n  2\n);\n\n// 2) \"Overview - ...\" (English objectives) \u2014 same split, English text.\nconst pOverviewEn = findParagraph(\"Overview - Introduce and teach concepts\");\nawait insertBreaksAfter(\n  pOverviewEn,\n  [\"and realize the changes that are necessary.\"],\n  2\n);\n\n// 3) \"Propriedade gerais ...\" (Portuguese short program) \u2014 single <w:br/>\n//    between \"estruturais\" and \"como estereoqu\u00edmica\".\nconst pPropriedade = findParagraph(\"Propriedade gerais dos compostos org\u00e2nicos\");\nawait insertBreaksAfter(\n  pPropriedade,\n  [\"\u00e1lcoois e caracter\u00edsticas estruturais\"],\n  1\n);\n\n// 4) \"1.Teoria de Bronsted ...\" (Portuguese numbered syllabus) \u2014 one <w:br/>\n//    before each numbered item (2. through 8.).\nconst pProgramaPt = findParagraph(\"1.Teoria de Bronsted e de Lewis\");\nawait insertBreaksAfter(\n  pProgramaPt,\n  [\n    \"acidez de compostos org\u00e2nicos\",\n    \"Rea\u00e7\u00e3o de Substitui\u00e7\u00e3o Radicalar. \",\n    \"Resolu\u00e7\u00e3o de Is\u00f4meros Espaciais.\",\n    \"Substitui\u00e7\u00e3o Nucleof\u00edlica, SN1, SN2, E1, E2. \",\n    \"Adi\u00e7\u00e3o conjugada em dienos (produto termodin\u00e2mico e cin\u00e9tico) \",\n    \"Infra-vermelho, Ultra-violeta e Fluorescencia \",\n    \"Rea\u00e7\u00e3o de Substitui\u00e7\u00e3o Nucleof\u00edlica.\"\n  ],\n  1\n);\n\n// 5) \"1.Bronsted and Lewis acid ...\" (English numbered syllabus) \u2014 same\n//    split structure as (4).\nconst pProgramaEn = findParagraph(\"1.Bronsted and Lewis acid of the organic compounds\");\nawait insertBreaksAfter(\n  pProgramaEn,\n  [\n    \"acid of the organic compounds\",\n    \"Radical substitution reaction.\",\n    \"resolution of stereoisomers.\",\n    \"Nucleophilic Substitution, SN1, SN2, E1, E2.\",\n    \"Conjugated Addition in dienes (thermodynamic and kinetic product).\",\n    \"UV and Fluorescence techniques\",\n    \"Aromatic Nucleophilic Substitution.\"\n  ],\n  1\n);\n\n// 6) \"Duas provas te\u00f3ricas ...\" (evaluation method) \u2014 blank line (two\n//    <w:br/>) between the first sentence and \"Aos alunos...\".\nconst pMetodo = findParagraph(\"Duas provas te\u00f3ricas e ao longo do semestre letivo\");\nawait insertBreaksAfter(\n  pMetodo,\n  [\"Duas provas te\u00f3ricas e ao longo do semestre letivo\"],\n  2\n);\n", "ps1": "# Replace long run-on paragraphs with explicit line-break (<w:br/>) separated\n# segments, matching the target OOXML diff. A Word manual line break\n# (Shift+Enter) is represented by Chr(11) (vertical tab); inserting that\n# character via Range.InsertAfter causes Word to store it as <w:br/> and\n# split the surrounding text into separate <w:t> runs - exactly what the\n# diff shows.\n\n$d = $word.ActiveDocument\n$c = [char]11\n\nfunction InsertBreaksAfter($para, $anchor, $count) {\n    # Finds $anchor inside $para's range (must be unique within it) and\n    # inserts $count manual line breaks immediately after the match,\n    # without touching any other text in the paragraph.\n    $r = $para.Range.Duplicate\n    $found = $r.Find.Execute($anchor)\n    if (-not $found) {\n        throw \"anchor not found: $anchor\"\n    }\n    $r.Collapse(0)  # wdCollapseEnd\n    $breaks = \"\"\n    for ($i = 0; $i -lt $count; $i++) {\n        $breaks = $breaks + $c\n    }\n    $r.InsertAfter($breaks)\n}\n\n# 1) \"Gerais - ...\" (Portuguese objectives) - split into two sentences with a\n#    blank line (two <w:br/>) between them.\nInsertBreaksAfter $d.Paragraphs(6) \"realizarem as mudan\u00e7as que se fa\u00e7am necess\u00e1rias.\" 2\n\n# 2) \"Overview - ...\" (English objectives) - same split, English text.\nInsertBreaksAfter $d.Paragraphs(7) \"and realize the changes that are necessary.\" 2\n\n# 3) \"Propriedade gerais ...\" (Portuguese short program) - single <w:br/>\n#    between \"estruturais\" and \"como estereoqu\u00edmica\".\nInsertBreaksAfter $d.Paragraphs(11) \"\u00e1lcoois e caracter\u00edsticas estruturais\" 1\n\n# 4) \"1.Teoria de Bronsted ...\" (Portuguese numbered syllabus) - one <w:br/>\n#    before each numbered item (2. through 8.).\nInsertBreaksAfter $d.Paragraphs(14) \"acidez de compostos org\u00e2nicos\" 1\nInsertBreaksAfter $d.Paragraphs(14) \"Rea\u00e7\u00e3o de Substitui\u00e7\u00e3o Radicalar. \" 1\nInsertBreaksAfter $d.Paragraphs(14) \"Resolu\u00e7\u00e3o de Is\u00f4meros Espaciais.\" 1\nInsertBreaksAfter $d.Paragraphs(14) \"Substitui\u00e7\u00e3o Nucleof\u00edlica, SN1, SN2, E1, E2. \" 1\nInsertBreaksAfter $d.Paragraphs(14) \"Adi\u00e7\u00e3o conjugada em dienos (produto termodin\u00e2mico e cin\u00e9tico) \" 1\nInsertBreaksAfter $d.Paragraphs(14) \"Infra-vermelho, Ultra-violeta e Fluorescencia \" 1\nInsertBreaksAfter $d.Paragraphs(14) \"Rea\u00e7\u00e3o de Substitui\u00e7\u00e3o Nucleof\u00edlica.\" 1\n\n# 5) \"1.Bronsted and Lewis acid ...\" (English numbered syllabus) - same split\n#    structure as (4).\nInsertBreaksAfter $d.Paragraphs(15) \"acid of the organic compounds\" 1\nInsertBreaksAfter $d.Paragraphs(15) \"Radical substitution reaction.\" 1\nInsertBreaksAfter $d.Paragraphs(15) \"resolution of stereoisomers.\" 1\nInsertBreaksAfter $d.Paragraphs(15) \"Nucleophilic Substitution, SN1, SN2, E1, E2.\" 1\nInsertBreaksAfter $d.Paragraphs(15) \"Conjugated Addition in dienes (thermodynamic and kinetic product).\" 1\nInsertBreaksAfter $d.Paragraphs(15) \"UV and Fluorescence techniques\" 1\nInsertBreaksAfter $d.Paragraphs(15) \"Aromatic Nucleophilic Substitution.\" 1\n\n# 6) \"Duas provas te\u00f3ricas ...\" (evaluation method) - blank line (two\n#    <w:br/>) between the first sentence and \"Aos alunos...\".\nInsertBreaksAfter $d.Paragraphs(17) \"Duas provas te\u00f3ricas e ao longo do semestre letivo\" 2\n"}
